$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Range("C6").Value = "Continued learning and practising Dart, Flutter"
